# Auto-generated Excel COM-interop PowerShell script
# Applies numeric cell updates to the Seraph_Profits workbook sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 15789
$ws.Range("J10").Value = 15789
$ws.Range("L10").Value = 15789
$ws.Range("N10").Value = -16375
$ws.Range("H33").Value = 250.22223
$ws.Range("I33").Value = 116.1875
$ws.Range("K33").Value = 116.1875
$ws.Range("M33").Value = 112.8125
$ws.Range("H40").Value = 1945.4546
$ws.Range("I40").Value = 1988.8889
$ws.Range("J40").Value = 1750
$ws.Range("K40").Value = 1988.8889
$ws.Range("L40").Value = 1750
$ws.Range("M40").Value = -1813.8889
$ws.Range("N40").Value = -2100
$ws.Range("H62").Value = 6706.364
$ws.Range("I62").Value = 5439.5
$ws.Range("J62").Value = 6987.8887
$ws.Range("K62").Value = 5439.5
$ws.Range("L62").Value = 6987.8887
$ws.Range("M62").Value = -4815.5
$ws.Range("N62").Value = -8235.8887
$ws.Range("H65").Value = 6706.364
$ws.Range("I65").Value = 5439.5
$ws.Range("J65").Value = 6987.8887
$ws.Range("K65").Value = 27197.5
$ws.Range("L65").Value = 34939.4435
$ws.Range("M65").Value = -24077.5
$ws.Range("N65").Value = -41179.4435
$ws.Range("H113").Value = 4101.4443
$ws.Range("I113").Value = 4244.7144
$ws.Range("J113").Value = 3600
$ws.Range("K113").Value = 4244.7144
$ws.Range("L113").Value = 3600
$ws.Range("M113").Value = -990.7143999999998
$ws.Range("N113").Value = -10108
$ws.Range("H137").Value = 3528.6072
$ws.Range("I137").Value = 1599.9375
$ws.Range("K137").Value = 4799.8125
$ws.Range("M137").Value = -2249.8125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 749.25
$ws.Range("I2").Value = 749.25
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 749.25
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -636.25
$ws.Range("N2").ClearContents()
$ws.Range("H45").Value = 2249.75
$ws.Range("J45").Value = 2249.75
$ws.Range("L45").Value = 2249.75
$ws.Range("N45").Value = -3003.75
$ws.Range("H102").Value = 1009.5
$ws.Range("I102").Value = 1010.8571
$ws.Range("K102").Value = 1010.8571
$ws.Range("M102").Value = 611.1429000000001
$ws.Range("H110").Value = 6710.6665
$ws.Range("I110").Value = 7771
$ws.Range("K110").Value = 7771
$ws.Range("M110").Value = -5726
$ws.Range("H116").Value = 749.25
$ws.Range("I116").Value = 749.25
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 749.25
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1544.75
$ws.Range("N116").ClearContents()
$ws.Range("H122").Value = 3004.9524
$ws.Range("J122").Value = 4708.857
$ws.Range("L122").Value = 14126.571
$ws.Range("N122").Value = -19026.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 749.25
$ws.Range("I3").Value = 749.25
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 749.25
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -635.25
$ws.Range("N3").ClearContents()
$ws.Range("H86").Value = 2631.2
$ws.Range("I86").Value = 1852.8334
$ws.Range("K86").Value = 1852.8334
$ws.Range("M86").Value = -729.8334
$ws.Range("H89").Value = 2631.2
$ws.Range("I89").Value = 1852.8334
$ws.Range("K89").Value = 9264.166999999999
$ws.Range("M89").Value = -3648.166999999999
$ws.Range("H99").Value = 1222.8572
$ws.Range("I99").Value = 1093
$ws.Range("K99").Value = 1093
$ws.Range("M99").Value = 405
$ws.Range("H134").Value = 2154.2173
$ws.Range("I134").Value = 1627.75
$ws.Range("J134").Value = 5664
$ws.Range("K134").Value = 4883.25
$ws.Range("L134").Value = 16992
$ws.Range("M134").Value = -2348.25
$ws.Range("N134").Value = -22062

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 135.45454
$ws.Range("I7").Value = 163.46666
$ws.Range("K7").Value = 163.46666
$ws.Range("M7").Value = -50.46665999999999
$ws.Range("H10").Value = 1016.6667
$ws.Range("I10").Value = 1016.6667
$ws.Range("K10").Value = 1016.6667
$ws.Range("M10").Value = -877.6667
$ws.Range("H122").Value = 916.5714
$ws.Range("I122").Value = 883.2
$ws.Range("K122").Value = 2649.6
$ws.Range("M122").Value = -199.6000000000004
$ws.Range("H132").Value = 3597.5
$ws.Range("I132").Value = 2522.5
$ws.Range("J132").Value = 5747.5
$ws.Range("K132").Value = 7567.5
$ws.Range("L132").Value = 17242.5
$ws.Range("M132").Value = -5037.5
$ws.Range("N132").Value = -22302.5
$ws.Range("H134").Value = 2409.16
$ws.Range("I134").Value = 1649.0454
$ws.Range("K134").Value = 4947.1362
$ws.Range("M134").Value = -2412.1362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 152
$ws.Range("I2").Value = 31.555555
$ws.Range("K2").Value = 189.33333
$ws.Range("M2").Value = -76.33332999999999
$ws.Range("H116").Value = 800
$ws.Range("I116").Value = 300
$ws.Range("K116").Value = 900
$ws.Range("M116").Value = 2542
$ws.Range("H122").Value = 1249
$ws.Range("I122").Value = 499
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 4491
$ws.Range("L122").Value = 17991
$ws.Range("M122").Value = -2041
$ws.Range("N122").Value = -22891
$ws.Range("H133").Value = 20011
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("H140").Value = 3927.3076
$ws.Range("I140").Value = 3005.0908
$ws.Range("K140").Value = 9015.2724
$ws.Range("M140").Value = -3835.2724
$ws.Range("H141").Value = 9199.200000000001
$ws.Range("I141").Value = 1999
$ws.Range("K141").Value = 5997
$ws.Range("M141").Value = -817

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 2000
$ws.Range("J14").Value = 2000
$ws.Range("L14").Value = 2000
$ws.Range("N14").Value = -2336
$ws.Range("H97").Value = 1246.8518
$ws.Range("I97").Value = 1320.9546
$ws.Range("K97").Value = 1320.9546
$ws.Range("M97").Value = -824.9546
$ws.Range("H102").Value = 2520.1333
$ws.Range("I102").Value = 1281.625
$ws.Range("J102").Value = 3935.5715
$ws.Range("K102").Value = 1281.625
$ws.Range("L102").Value = 3935.5715
$ws.Range("M102").Value = 340.375
$ws.Range("N102").Value = -7179.5715
$ws.Range("H122").Value = 422371.1
$ws.Range("I122").Value = 4239
$ws.Range("K122").Value = 12717
$ws.Range("M122").Value = -10267

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1291.6
$ws.Range("I22").Value = 643
$ws.Range("J22").Value = 1940.2
$ws.Range("K22").Value = 643
$ws.Range("L22").Value = 1940.2
$ws.Range("M22").Value = -348
$ws.Range("N22").Value = -2530.2
$ws.Range("H27").Value = 1291.6
$ws.Range("I27").Value = 643
$ws.Range("J27").Value = 1940.2
$ws.Range("K27").Value = 643
$ws.Range("L27").Value = 1940.2
$ws.Range("M27").Value = -536
$ws.Range("N27").Value = -2154.2
$ws.Range("H40").Value = 3970.4285
$ws.Range("I40").Value = 3821.111
$ws.Range("K40").Value = 3821.111
$ws.Range("M40").Value = -3685.111
$ws.Range("H61").Value = 5491.6875
$ws.Range("I61").Value = 7183.2856
$ws.Range("J61").Value = 4176
$ws.Range("K61").Value = 7183.2856
$ws.Range("L61").Value = 4176
$ws.Range("M61").Value = -6981.2856
$ws.Range("N61").Value = -4580
$ws.Range("H113").Value = 5491.6875
$ws.Range("I113").Value = 7183.2856
$ws.Range("J113").Value = 4176
$ws.Range("K113").Value = 7183.2856
$ws.Range("L113").Value = 4176
$ws.Range("M113").Value = -5013.2856
$ws.Range("N113").Value = -8516
$ws.Range("H132").Value = 6459.3
$ws.Range("I132").Value = 6398.5
$ws.Range("J132").Value = 6499.8335
$ws.Range("K132").Value = 19195.5
$ws.Range("L132").Value = 19499.5005
$ws.Range("M132").Value = -16665.5
$ws.Range("N132").Value = -24559.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 2562.25
$ws.Range("I6").Value = 2125
$ws.Range("J6").Value = 2999.5
$ws.Range("K6").Value = 2125
$ws.Range("L6").Value = 2999.5
$ws.Range("M6").Value = -2010
$ws.Range("N6").Value = -3229.5
$ws.Range("H7").Value = 3621.8
$ws.Range("I7").Value = 8002
$ws.Range("K7").Value = 8002
$ws.Range("M7").Value = -7889
$ws.Range("H132").Value = 2764.9473
$ws.Range("I132").Value = 1937
$ws.Range("K132").Value = 5811
$ws.Range("M132").Value = -2764.9473
